$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 10: add E10 = "?" and change F10 value
$ws.Range("E10").Value = "?"
$ws.Range("F10").Value = 16118

# Row 11: add E11 = "?", add F11 formula (with F10's number format/border style), change G11 formula
$ws.Range("E11").Value = "?"
$ws.Range("F10").Copy()
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F11").Formula = "=49525-17058-16118"
$ws.Range("G11").Formula = "=G10+F10"

# Row 12: add E12 = "?", add F12 value (with F10's style), change G12 formula
$ws.Range("E12").Value = "?"
$ws.Range("F10").Copy()
$ws.Range("F12").PasteSpecial(-4122)
$ws.Range("F12").Value = 17058
$ws.Range("G12").Formula = "=G11+F11"

# Row 13: change F13 formula
$ws.Range("F13").Formula = "=51590-17687-17232"

# Row 14: add F14 value
$ws.Range("F14").Value = 17232

# Update the selection to F10 as in the diff
$ws.Range("F10").Select()
